$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The per-jefe support email addresses duplicated the general support
# email already used elsewhere in the sheet, so replace them with the
# shared address and omit the now-redundant strings.
$email = "mariavyeguezp@gmail.com"

$ws.Range("C2").Value = $email
$ws.Range("C3").Value = $email
$ws.Range("C4").Value = $email

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:" + $email) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3:C4"), "mailto:" + $email, [Type]::Missing, [Type]::Missing, $email) | Out-Null

# Adding the hyperlinks applies direct font formatting; restore the
# shared "Hipervínculo" cell style used by the other hyperlink cells.
$ws.Range("C2:C4").Style = "Hipervínculo"

# Update the saved view: clear the scrolled-away top-left cell and move
# the active selection.
$ws.Range("E14").Select()
